$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text before writing, so numeric-looking strings
# (e.g. "545.38") are preserved exactly as text instead of being coerced
# into floating-point numbers by Excel.
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '60.356.09'
$ws.Cells.Item(2, 5).Value = '  +0.28%  '

$ws.Cells.Item(3, 4).Value = '2.317.83'
$ws.Cells.Item(3, 5).Value = '  -1.40%  '

$ws.Cells.Item(4, 5).Value = '  -0.13%  '

$ws.Cells.Item(5, 4).Value = '545.38'
$ws.Cells.Item(5, 5).Value = '  -0.03%  '

$ws.Cells.Item(6, 4).Value = '129.75'
$ws.Cells.Item(6, 5).Value = '  -1.96%  '

$ws.Cells.Item(7, 5).Value = '  -0.09%  '

$ws.Cells.Item(8, 5).Value = '  -2.29%  '

$ws.Cells.Item(9, 4).Value = '2.318.24'
$ws.Cells.Item(9, 5).Value = '  -1.25%  '

$ws.Cells.Item(10, 5).Value = '  -0.14%  '

$ws.Cells.Item(11, 4).Value = '5.54'
$ws.Cells.Item(11, 5).Value = '  +0.42%  '

$ws.Cells.Item(12, 5).Value = '  -0.53%  '

$ws.Cells.Item(13, 4).Value = '0.335'
$ws.Cells.Item(13, 5).Value = '  +0.30%  '

$ws.Cells.Item(14, 4).Value = '23.45'
$ws.Cells.Item(14, 5).Value = '  -2.00%  '

$ws.Cells.Item(15, 4).Value = '60.370.24'
$ws.Cells.Item(15, 5).Value = '  +0.18%  '

$ws.Cells.Item(16, 4).Value = '2.730.09'
$ws.Cells.Item(16, 5).Value = '  -1.54%  '

$ws.Cells.Item(17, 5).Value = '  -0.01%  '

$ws.Cells.Item(18, 4).Value = '2.321.56'
$ws.Cells.Item(18, 5).Value = '  -1.25%  '

$ws.Cells.Item(19, 4).Value = '10.57'
$ws.Cells.Item(19, 5).Value = '  -1.18%  '

$ws.Cells.Item(20, 5).Value = '  -2.58%  '

$ws.Cells.Item(21, 4).Value = '313.09'
$ws.Cells.Item(21, 5).Value = '  -0.56%  '

$ws.Cells.Item(22, 5).Value = '  -4.03%  '

$ws.Cells.Item(23, 4).Value = '1.00'
$ws.Cells.Item(23, 5).Value = '  -0.01%  '

$ws.Cells.Item(24, 4).Value = '64.04'
$ws.Cells.Item(24, 5).Value = '  +1.30%  '

$ws.Cells.Item(25, 5).Value = '  -0.22%  '

$ws.Cells.Item(26, 4).Value = '0.999'
$ws.Cells.Item(26, 5).Value = '  -0.23%  '

$ws.Cells.Item(27, 4).Value = '7.88'
$ws.Cells.Item(27, 5).Value = '  -0.89%  '

$ws.Cells.Item(28, 5).Value = '  +2.54%  '

$ws.Cells.Item(29, 5).Value = '  +6.92%  '

$ws.Cells.Item(30, 4).Value = '173.20'
$ws.Cells.Item(30, 5).Value = '  +1.09%  '

$ws.Cells.Item(31, 5).Value = '  -1.53%  '

$ws.Cells.Item(32, 5).Value = '  -0.67%  '

$ws.Cells.Item(33, 4).Value = '5.96'
$ws.Cells.Item(33, 5).Value = '  +0.30%  '

$ws.Cells.Item(34, 5).Value = '  +0.00%  '

$ws.Cells.Item(35, 5).Value = '  -3.91%  '

$ws.Cells.Item(36, 5).Value = '  -0.30%  '

$ws.Cells.Item(38, 5).Value = '  -0.08%  '

$ws.Cells.Item(39, 5).Value = '  -1.12%  '

$ws.Cells.Item(40, 4).Value = '317.37'
$ws.Cells.Item(40, 5).Value = '  -0.03%  '

$ws.Cells.Item(41, 4).Value = '38.02'
$ws.Cells.Item(41, 5).Value = '  -0.46%  '

$ws.Cells.Item(42, 5).Value = '  -0.58%  '

$ws.Cells.Item(43, 4).Value = '136.83'
$ws.Cells.Item(43, 5).Value = '  -4.07%  '

$ws.Cells.Item(44, 5).Value = '  +0.99%  '

$ws.Cells.Item(46, 4).Value = '19.07'
$ws.Cells.Item(46, 5).Value = '  -1.16%  '

$ws.Cells.Item(47, 4).Value = '0.565'
$ws.Cells.Item(47, 5).Value = '  +0.35%  '

$ws.Cells.Item(48, 4).Value = '0.0495'
$ws.Cells.Item(48, 5).Value = '  -0.67%  '

$ws.Cells.Item(49, 5).Value = '  +0.46%  '

$ws.Cells.Item(50, 4).Value = '0.0₆0218'
$ws.Cells.Item(50, 5).Value = '  +6.24%  '

$ws.Cells.Item(51, 4).Value = '10.96'
$ws.Cells.Item(51, 5).Value = '  -0.77%  '

# Restore the original (default) style on column D so no stray
# number-format styling is left behind on the cells.
$colD.Style = "Normal"
